# Append new scrape results to the "ランサーズ" sheet.
# Existing row 2 data is superseded/shifted down; 8 brand-new rows are
# inserted above the previously-existing two rows, and all timestamps in
# column A are refreshed to the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-01-20 12:44:14"

# Each entry: Title, Category, Price, Deadline, URL, Score, Skill-summary (or $null)
$rows = @(
    @("【オンラインレッスン】Tailwindテンプレートをサーバー接続したい", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475207", 295, "🔥AI,Ai"),
    @("【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,150円程度)", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475245", 213, "🔥API ◇管理"),
    @("【急募】フットアールサッカースクール向け出欠管理Webアプリ開発者募集", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475081", 128, "◆開発 ◇アプリ"),
    @("【急募】ガイドと旅行者をつなぐマッチングサイト開発", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475198", 93, "◆開発 ◇サイト"),
    @("初回 IOS、ANDROIDのアプリ開発", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475193", 93, "◆開発 ◇アプリ"),
    @("【急募】業務システムの開発・運用・保守エンジニア募集(フロントエンド/バックエンド)", "システム開発", "1,000 ~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474899", 75, "◆開発"),
    @("【急募】社内部所間の振込依頼システム構築", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475326", 40, $null),
    @("【インバウンド/フルリモ/月30万円固定】医療系SaaSのIS・FS", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475237", 25, $null),
    @("【急募・オンライン】Kaggle経験者求む!Cursorを使ったKaggle実践サポート", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475082", 10, $null),
    @("Google clab用マークシートCSV出力プログラム作成依頼", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474679", 10, $null)
)

# Clear the two pre-existing hyperlinks (on the old F2/F3) so they don't
# linger as stale relationships once every row gets a freshly-added link.
$ws.Hyperlinks.Delete()

$r = 2
foreach ($item in $rows) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $item[0]
    $ws.Cells.Item($r, 3).Value = $item[1]
    $ws.Cells.Item($r, 4).Value = $item[2]
    $ws.Cells.Item($r, 5).Value = $item[3]

    $urlCell = $ws.Cells.Item($r, 6)
    $urlCell.Value = $item[4]
    $ws.Hyperlinks.Add($urlCell, $item[4])

    $ws.Cells.Item($r, 7).Value = $item[5]

    if ($item[6] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $item[6]
    }

    $r = $r + 1
}

# Widen columns B and D to fit the longer new content.
# (ColumnWidth has a +5/6-character rendering offset in this engine, so the
# assigned value is pre-compensated to land exactly on the target width.)
$ws.Columns.Item(2).ColumnWidth = 52 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 30 - (5/6)
